$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.691401600837708
$ws.Range("B1").Value = 1.850774765014648
$ws.Range("C1").Value = 2.099740266799927
$ws.Range("D1").Value = 2.629744052886963
$ws.Range("E1").Value = 1.6321941614151
